# Generate Report for Handoff
#
# This updates the "65f09445-3d6a-40c8-8e93-e3c21a26c086" row on every sheet
# so that it reflects that the file is now "Ready for handoff" (instead of
# "Handed back: in sync with en-US"), refreshes the handback timestamps, and
# records the handback-out-of-date error message + widens the Error Detail
# column so the long message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/925d297afd8fa64c08f3e77058ccfe3047070021/e2e/65f09445-3d6a-40c8-8e93-e3c21a26c086.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec10d16887e68ca8af6ff84f888cd96dc28f6840/e2e/65f09445-3d6a-40c8-8e93-e3c21a26c086.md."

# --- Overview sheet: row 3 is the 65f09445-...md file ---------------------
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-06 09:05:45"

# --- zh-cn sheet: row 3 is the 65f09445-...md file -------------------------
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-06 09:05:39"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").EntireColumn.ColumnWidth = 39.17

# --- de-de sheet: row 3 is the 65f09445-...md file --------------------------
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-06 09:05:45"
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").EntireColumn.ColumnWidth = 39.17
